# Update countries & provincias Spain
# Applies the 1-Aug-2020 11:41 data refresh to the "Pais" sheet:
#   - Updates the "last updated" timestamp banner
#   - Re-ranks two pairs of countries whose case counts crossed over
#     (Filipinas now ahead of Egipto; Hong Kong now ahead of Tailandia)
#   - Updates the statistic columns (B..H) for every country whose
#     figures changed in this refresh

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 1 de Agosto de 2020 a las 11:41"

# --- Re-rank: Indonesia / Egipto / Filipinas block (rows 27-29) -------
# Filipinas' total now exceeds Egipto's (which is unchanged), so the two
# swap places; Egipto's figures carry over untouched to its new row.
$ws.Range("A28").Value = "Filipinas"
$ws.Range("A29").Value = "Egipto"

# --- Re-rank: Libia / Tailandia / Hong Kong block (rows 110-113) ------
# Hong Kong's total now exceeds Tailandia's (which is unchanged), so the
# two swap places; Tailandia's figures carry over untouched to its new row.
$ws.Range("A111").Value = "Hong Kong"
$ws.Range("A112").Value = "Tailandia"

# --- Statistic refresh (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---------------

# Estados Unidos
$ws.Range("B4").Value = 4706180
$ws.Range("C4").Value = 291
$ws.Range("D4").Value = 2328445
$ws.Range("E4").Value = 2220971
$ws.Range("G4").Value = 17
$ws.Range("H4").Value = 156764

# Banglades
$ws.Range("B19").Value = 239860
$ws.Range("C19").Value = 2199
$ws.Range("D19").Value = 136253
$ws.Range("E19").Value = 100475
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = 3132

# Indonesia
$ws.Range("B27").Value = 109936
$ws.Range("C27").Value = 1560
$ws.Range("D27").Value = 67919
$ws.Range("E27").Value = 36824
$ws.Range("G27").Value = 62
$ws.Range("H27").Value = 5193

# Filipinas (now row 28)
$ws.Range("B28").Value = 98232
$ws.Range("C28").Value = 4963
$ws.Range("D28").Value = 65265
$ws.Range("E28").Value = 30928
$ws.Range("G28").Value = 17
$ws.Range("H28").Value = 2039

# Egipto (now row 29, figures unchanged from previous refresh)
$ws.Range("B29").Value = 94078
$ws.Range("C29").Value = 0
$ws.Range("D29").Value = 39638
$ws.Range("E29").Value = 49635
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 4805

# Polonia
$ws.Range("B49").Value = 46346
$ws.Range("C49").Value = 658
$ws.Range("D49").Value = 34374
$ws.Range("E49").Value = 10251
$ws.Range("G49").Value = 5
$ws.Range("H49").Value = 1721

# Moldavia
$ws.Range("D63").Value = 17571
$ws.Range("E63").Value = 6379
$ws.Range("G63").Value = 5
$ws.Range("H63").Value = 783

# Austria
$ws.Range("B66").Value = 21212
$ws.Range("C66").Value = 82
$ws.Range("D66").Value = 18911
$ws.Range("E66").Value = 1583

# Consejo Danes para los Refugiados
$ws.Range("B87").Value = 9084
$ws.Range("C87").Value = 14
$ws.Range("D87").Value = 7030
$ws.Range("E87").Value = 1839

# Malasia
$ws.Range("B88").Value = 8985
$ws.Range("C88").Value = 9
$ws.Range("D88").Value = 8647
$ws.Range("E88").Value = 213

# Finlandia
$ws.Range("B90").Value = 7443
$ws.Range("C90").Value = 11
$ws.Range("E90").Value = 164

# Hong Kong (now row 111)
$ws.Range("B111").Value = 3398
$ws.Range("C111").Value = 125
$ws.Range("D111").Value = 1858
$ws.Range("E111").Value = 1509
$ws.Range("G111").Value = 4
$ws.Range("H111").Value = 31

# Tailandia (now row 112, figures unchanged from previous refresh)
$ws.Range("B112").Value = 3312
$ws.Range("C112").Value = 2
$ws.Range("D112").Value = 3135
$ws.Range("E112").Value = 119
$ws.Range("G112").Value = 0
$ws.Range("H112").Value = 58

# Sri Lanka
$ws.Range("D118").Value = 2439
$ws.Range("E118").Value = 365

# Eslovenia
$ws.Range("B125").Value = 2171
$ws.Range("C125").Value = 15
$ws.Range("D125").Value = 1821
$ws.Range("E125").Value = 231
